$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 499.5
$ws.Range("I7").Value = 499.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 499.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -387.5
$ws.Range("N7").ClearContents()

$ws.Range("H14").Value = 499.5
$ws.Range("I14").Value = 499.5
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 499.5
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -308.5
$ws.Range("N14").ClearContents()

$ws.Range("H17").Value = 956.25
$ws.Range("J17").Value = 956.25
$ws.Range("L17").Value = 2868.75
$ws.Range("N17").Value = -3204.75

$ws.Range("H33").Value = 671.46155
$ws.Range("I33").Value = 695.52
$ws.Range("K33").Value = 695.52
$ws.Range("M33").Value = -466.52

$ws.Range("H107").Value = 924.2857
$ws.Range("I107").Value = 839.2308
$ws.Range("J107").Value = 1062.5
$ws.Range("K107").Value = 839.2308
$ws.Range("L107").Value = 1062.5
$ws.Range("M107").Value = 1080.7692
$ws.Range("N107").Value = -4902.5

$ws.Range("H112").Value = 3049.2856
$ws.Range("J112").Value = 3143.7036
$ws.Range("L112").Value = 9431.110799999999
$ws.Range("N112").Value = -11647.1108

$ws.Range("H137").Value = 20236.547
$ws.Range("I137").Value = 2087.6667
$ws.Range("J137").Value = 32146.75
$ws.Range("K137").Value = 6263.000100000001
$ws.Range("L137").Value = 96440.25
$ws.Range("M137").Value = -3713.000100000001
$ws.Range("N137").Value = -101540.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 808.0476
$ws.Range("I2").Value = 745.2941
$ws.Range("J2").Value = 1074.75
$ws.Range("K2").Value = 745.2941
$ws.Range("L2").Value = 1074.75
$ws.Range("M2").Value = -632.2941
$ws.Range("N2").Value = -1300.75

$ws.Range("H32").Value = 6834.98
$ws.Range("I32").Value = 5864.109
$ws.Range("J32").Value = 18000
$ws.Range("K32").Value = 5864.109
$ws.Range("L32").Value = 18000
$ws.Range("M32").Value = -5577.109
$ws.Range("N32").Value = -18574

$ws.Range("H45").Value = 1510.125
$ws.Range("I45").Value = 1212.3334
$ws.Range("J45").Value = 2403.5
$ws.Range("K45").Value = 1212.3334
$ws.Range("L45").Value = 2403.5
$ws.Range("M45").Value = -835.3334
$ws.Range("N45").Value = -3157.5

$ws.Range("H116").Value = 808.0476
$ws.Range("I116").Value = 745.2941
$ws.Range("J116").Value = 1074.75
$ws.Range("K116").Value = 745.2941
$ws.Range("L116").Value = 1074.75
$ws.Range("M116").Value = 1548.7059
$ws.Range("N116").Value = -5662.75

$ws.Range("H122").Value = 928.4783
$ws.Range("I122").Value = 892.3684
$ws.Range("J122").Value = 1100
$ws.Range("K122").Value = 2677.1052
$ws.Range("L122").Value = 3300
$ws.Range("M122").Value = -227.1052
$ws.Range("N122").Value = -8200

$ws.Range("H132").Value = 725600.2
$ws.Range("I132").Value = 1325579.8
$ws.Range("J132").Value = 5624.6333
$ws.Range("K132").Value = 3976739.4
$ws.Range("L132").Value = 16873.8999
$ws.Range("M132").Value = -3974209.4
$ws.Range("N132").Value = -21933.8999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 808.0476
$ws.Range("I3").Value = 745.2941
$ws.Range("J3").Value = 1074.75
$ws.Range("K3").Value = 745.2941
$ws.Range("L3").Value = 1074.75
$ws.Range("M3").Value = -631.2941
$ws.Range("N3").Value = -1302.75

$ws.Range("H112").Value = 22250
$ws.Range("J112").Value = 22250
$ws.Range("L112").Value = 22250
$ws.Range("N112").Value = -25204

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1257.8572
$ws.Range("I107").Value = 287.14285
$ws.Range("J107").Value = 2228.5715
$ws.Range("K107").Value = 287.14285
$ws.Range("L107").Value = 2228.5715
$ws.Range("M107").Value = 1632.85715
$ws.Range("N107").Value = -6068.5715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 1501
$ws.Range("I19").Value = 1501
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 4503
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -4329
$ws.Range("N19").ClearContents()

$ws.Range("H113").Value = 433.5625
$ws.Range("I113").Value = 349.41666
$ws.Range("J113").Value = 686
$ws.Range("K113").Value = 1048.24998
$ws.Range("L113").Value = 2058
$ws.Range("M113").Value = 1121.75002
$ws.Range("N113").Value = -6398

$ws.Range("H117").Value = 14291317
$ws.Range("I117").Value = 10029
$ws.Range("J117").Value = 15878126
$ws.Range("K117").Value = 30087
$ws.Range("L117").Value = 47634378
$ws.Range("M117").Value = -26645
$ws.Range("N117").Value = -47641262

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2372.4375
$ws.Range("I113").Value = 1655.6666
$ws.Range("J113").Value = 3294
$ws.Range("K113").Value = 1655.6666
$ws.Range("L113").Value = 3294
$ws.Range("M113").Value = 514.3334
$ws.Range("N113").Value = -7634

$ws.Range("H132").Value = 2256.75
$ws.Range("I132").Value = 1456.9333
$ws.Range("J132").Value = 3589.7778
$ws.Range("K132").Value = 4370.7999
$ws.Range("L132").Value = 10769.3334
$ws.Range("M132").Value = -1840.7999
$ws.Range("N132").Value = -15829.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 142860060
$ws.Range("I40").Value = 200001500
$ws.Range("J40").Value = 6490
$ws.Range("K40").Value = 200001500
$ws.Range("L40").Value = 6490
$ws.Range("M40").Value = -200001364
$ws.Range("N40").Value = -6762

$ws.Range("H68").Value = 2211.25
$ws.Range("I68").Value = 1800
$ws.Range("J68").Value = 3993.3333
$ws.Range("K68").Value = 1800
$ws.Range("L68").Value = 3993.3333
$ws.Range("M68").Value = -1051
$ws.Range("N68").Value = -5491.3333

$ws.Range("H71").Value = 2211.25
$ws.Range("I71").Value = 1800
$ws.Range("J71").Value = 3993.3333
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 19966.6665
$ws.Range("M71").Value = -5256
$ws.Range("N71").Value = -27454.6665

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H100").Value = 4608.9
$ws.Range("I100").Value = 5500
$ws.Range("J100").Value = 4227
$ws.Range("K100").Value = 5500
$ws.Range("L100").Value = 4227
$ws.Range("M100").Value = -4959
$ws.Range("N100").Value = -5309

$ws.Range("H122").Value = 9544.294
$ws.Range("I122").Value = 13967.111
$ws.Range("J122").Value = 4568.625
$ws.Range("K122").Value = 41901.333
$ws.Range("L122").Value = 13705.875
$ws.Range("M122").Value = -39451.333
$ws.Range("N122").Value = -18605.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1808.1842
$ws.Range("I132").Value = 1721.3572
$ws.Range("J132").Value = 2051.3
$ws.Range("K132").Value = 5164.071599999999
$ws.Range("L132").Value = 6153.900000000001
$ws.Range("M132").Value = -2634.071599999999
$ws.Range("N132").Value = -11213.9
